$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 1269006
$ws_ALC.Range("J17").Value = 1293406.8
$ws_ALC.Range("L17").Value = 3880220.4
$ws_ALC.Range("N17").Value = -3880556.4

$ws_ALC.Range("H33").Value = 9298
$ws_ALC.Range("I33").Value = 11505.368
$ws_ALC.Range("K33").Value = 11505.368
$ws_ALC.Range("M33").Value = -11276.368

$ws_ALC.Range("H39").Value = 987.5
$ws_ALC.Range("I39").Value = 59.46154
$ws_ALC.Range("J39").Value = 3400.4
$ws_ALC.Range("K39").Value = 178.38462
$ws_ALC.Range("L39").Value = 10201.2
$ws_ALC.Range("M39").Value = 117.61538
$ws_ALC.Range("N39").Value = -10793.2

$ws_ALC.Range("H51").Value = 4881.0527
$ws_ALC.Range("J51").Value = 4930
$ws_ALC.Range("L51").Value = 4930
$ws_ALC.Range("N51").Value = -5898

$ws_ALC.Range("H64").Value = 6476.222
$ws_ALC.Range("I64").Value = 5324.5
$ws_ALC.Range("J64").Value = 7397.6
$ws_ALC.Range("K64").Value = 5324.5
$ws_ALC.Range("L64").Value = 7397.6
$ws_ALC.Range("M64").Value = -5076.5
$ws_ALC.Range("N64").Value = -7893.6

$ws_ALC.Range("H67").Value = 6476.222
$ws_ALC.Range("I67").Value = 5324.5
$ws_ALC.Range("J67").Value = 7397.6
$ws_ALC.Range("K67").Value = 5324.5
$ws_ALC.Range("L67").Value = 7397.6
$ws_ALC.Range("M67").Value = -4466.5
$ws_ALC.Range("N67").Value = -9113.6

$ws_ALC.Range("H127").Value = 2866.875
$ws_ALC.Range("I127").Value = 766.6667
$ws_ALC.Range("J127").Value = 4127
$ws_ALC.Range("K127").Value = 2300.0001
$ws_ALC.Range("L127").Value = 12381
$ws_ALC.Range("M127").Value = 2659.9999
$ws_ALC.Range("N127").Value = -22301

$ws_ALC.Range("H141").Value = 2029.3
$ws_ALC.Range("I141").Value = 1569.1428
$ws_ALC.Range("K141").Value = 4707.428400000001
$ws_ALC.Range("M141").Value = 472.5715999999993

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 4476.533
$ws_ARM.Range("I2").Value = 4145.6
$ws_ARM.Range("J2").Value = 5138.4
$ws_ARM.Range("K2").Value = 4145.6
$ws_ARM.Range("L2").Value = 5138.4
$ws_ARM.Range("M2").Value = -4032.6
$ws_ARM.Range("N2").Value = -5364.4

$ws_ARM.Range("H32").Value = 4443.771
$ws_ARM.Range("I32").Value = 4434.159
$ws_ARM.Range("K32").Value = 4434.159
$ws_ARM.Range("M32").Value = -4147.159

$ws_ARM.Range("H42").Value = 5001
$ws_ARM.Range("I42").Value = 5001
$ws_ARM.Range("K42").Value = 5001
$ws_ARM.Range("M42").Value = -4515

$ws_ARM.Range("H61").Value = 3439.12
$ws_ARM.Range("I61").Value = 3189.7222
$ws_ARM.Range("K61").Value = 3189.7222
$ws_ARM.Range("M61").Value = -2977.7222

$ws_ARM.Range("H74").Value = 37203.97
$ws_ARM.Range("I74").Value = 40524.367
$ws_ARM.Range("J74").Value = 4000
$ws_ARM.Range("K74").Value = 40524.367
$ws_ARM.Range("L74").Value = 4000
$ws_ARM.Range("M74").Value = -39650.367
$ws_ARM.Range("N74").Value = -5748

$ws_ARM.Range("H77").Value = 37203.97
$ws_ARM.Range("I77").Value = 40524.367
$ws_ARM.Range("J77").Value = 4000
$ws_ARM.Range("K77").Value = 202621.835
$ws_ARM.Range("L77").Value = 20000
$ws_ARM.Range("M77").Value = -198253.835
$ws_ARM.Range("N77").Value = -28736

$ws_ARM.Range("H97").Value = 1131.5652
$ws_ARM.Range("I97").Value = 1019.8421
$ws_ARM.Range("J97").Value = 1662.25
$ws_ARM.Range("K97").Value = 1019.8421
$ws_ARM.Range("L97").Value = 1662.25
$ws_ARM.Range("M97").Value = -523.8421
$ws_ARM.Range("N97").Value = -2654.25

$ws_ARM.Range("H102").Value = 98728.48
$ws_ARM.Range("I102").Value = 108015.69
$ws_ARM.Range("K102").Value = 108015.69
$ws_ARM.Range("M102").Value = -106393.69

$ws_ARM.Range("H116").Value = 4476.533
$ws_ARM.Range("I116").Value = 4145.6
$ws_ARM.Range("J116").Value = 5138.4
$ws_ARM.Range("K116").Value = 4145.6
$ws_ARM.Range("L116").Value = 5138.4
$ws_ARM.Range("M116").Value = -1851.6
$ws_ARM.Range("N116").Value = -9726.4

$ws_ARM.Range("H132").Value = 72309.84
$ws_ARM.Range("I132").Value = 10114.728
$ws_ARM.Range("K132").Value = 30344.184
$ws_ARM.Range("M132").Value = -27814.184

$ws_ARM.Range("H136").Value = 3439.12
$ws_ARM.Range("I136").Value = 3189.7222
$ws_ARM.Range("K136").Value = 9569.1666
$ws_ARM.Range("M136").Value = -7019.1666

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 4476.533
$ws_BSM.Range("I3").Value = 4145.6
$ws_BSM.Range("J3").Value = 5138.4
$ws_BSM.Range("K3").Value = 4145.6
$ws_BSM.Range("L3").Value = 5138.4
$ws_BSM.Range("M3").Value = -4031.6
$ws_BSM.Range("N3").Value = -5366.4

$ws_BSM.Range("H44").Value = 28266.666
$ws_BSM.Range("I44").Value = 28266.666
$ws_BSM.Range("K44").Value = 28266.666
$ws_BSM.Range("M44").Value = -27769.666

$ws_BSM.Range("H105").Value = 14863.348
$ws_BSM.Range("I105").Value = 28910.875
$ws_BSM.Range("K105").Value = 28910.875
$ws_BSM.Range("M105").Value = -27163.875

$ws_BSM.Range("H132").Value = 139779.5
$ws_BSM.Range("J132").Value = 139779.5
$ws_BSM.Range("L132").Value = 139779.5
$ws_BSM.Range("N132").Value = -149899.5

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 3853.5
$ws_CRP.Range("I31").Value = 1565.6
$ws_CRP.Range("K31").Value = 1565.6
$ws_CRP.Range("M31").Value = -1270.6

$ws_CRP.Range("H34").Value = 3853.5
$ws_CRP.Range("I34").Value = 1565.6
$ws_CRP.Range("K34").Value = 1565.6
$ws_CRP.Range("M34").Value = -1363.6

$ws_CRP.Range("H58").Value = 18336844
$ws_CRP.Range("I58").Value = 4999.857
$ws_CRP.Range("J58").Value = 34377210
$ws_CRP.Range("K58").Value = 4999.857
$ws_CRP.Range("L58").Value = 34377210
$ws_CRP.Range("M58").Value = -4796.857
$ws_CRP.Range("N58").Value = -34377616

$ws_CRP.Range("H99").Value = 5643
$ws_CRP.Range("I99").Value = 5987.222
$ws_CRP.Range("K99").Value = 5987.222
$ws_CRP.Range("M99").Value = -4489.222

$ws_CRP.Range("H122").Value = 1640
$ws_CRP.Range("I122").Value = 1813.75
$ws_CRP.Range("J122").Value = 250
$ws_CRP.Range("K122").Value = 5441.25
$ws_CRP.Range("L122").Value = 750
$ws_CRP.Range("M122").Value = -2991.25
$ws_CRP.Range("N122").Value = -5650

$ws_CRP.Range("H126").Value = 5643
$ws_CRP.Range("I126").Value = 5987.222
$ws_CRP.Range("K126").Value = 17961.666
$ws_CRP.Range("M126").Value = -15491.666

$ws_CRP.Range("H132").Value = 44900.562
$ws_CRP.Range("I132").Value = 48476.727
$ws_CRP.Range("J132").Value = 5562.75
$ws_CRP.Range("K132").Value = 145430.181
$ws_CRP.Range("L132").Value = 16688.25
$ws_CRP.Range("M132").Value = -142900.181
$ws_CRP.Range("N132").Value = -21748.25

$ws_CRP.Range("H136").Value = 18336844
$ws_CRP.Range("I136").Value = 4999.857
$ws_CRP.Range("J136").Value = 34377210
$ws_CRP.Range("K136").Value = 14999.571
$ws_CRP.Range("L136").Value = 103131630
$ws_CRP.Range("M136").Value = -12449.571
$ws_CRP.Range("N136").Value = -103136730

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H34").Value = 577.4167
$ws_CUL.Range("I34").Value = 131.66667
$ws_CUL.Range("J34").Value = 1023.1667
$ws_CUL.Range("K34").Value = 395.00001
$ws_CUL.Range("L34").Value = 3069.5001
$ws_CUL.Range("M34").Value = -311.00001
$ws_CUL.Range("N34").Value = -3237.5001

$ws_CUL.Range("H102").Value = 2499
$ws_CUL.Range("I102").Value = 2499
$ws_CUL.Range("J102").Value = 0
$ws_CUL.Range("K102").Value = 7497
$ws_CUL.Range("L102").Value = 0
$ws_CUL.Range("M102").Value = -5063
$ws_CUL.Range("N102").ClearContents()

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H102").Value = 3770.2334
$ws_GSM.Range("J102").Value = 6602.4165
$ws_GSM.Range("L102").Value = 6602.4165
$ws_GSM.Range("N102").Value = -9846.416499999999

$ws_GSM.Range("H122").Value = 1649.5652
$ws_GSM.Range("I122").Value = 1180.5
$ws_GSM.Range("K122").Value = 3541.5
$ws_GSM.Range("M122").Value = -1091.5

$ws_GSM.Range("H126").Value = 3875.6667
$ws_GSM.Range("I126").Value = 3555.1428
$ws_GSM.Range("J126").Value = 4997.5
$ws_GSM.Range("K126").Value = 10665.4284
$ws_GSM.Range("L126").Value = 14992.5
$ws_GSM.Range("M126").Value = -8195.428400000001
$ws_GSM.Range("N126").Value = -19932.5

$ws_GSM.Range("H132").Value = 2600
$ws_GSM.Range("I132").Value = 2600
$ws_GSM.Range("K132").Value = 7800
$ws_GSM.Range("M132").Value = -5270

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 3225
$ws_LTW.Range("J22").Value = 4500
$ws_LTW.Range("L22").Value = 4500
$ws_LTW.Range("N22").Value = -5090

$ws_LTW.Range("H27").Value = 3225
$ws_LTW.Range("J27").Value = 4500
$ws_LTW.Range("L27").Value = 4500
$ws_LTW.Range("N27").Value = -4714

$ws_LTW.Range("H40").Value = 2055
$ws_LTW.Range("I40").Value = 1807.1666
$ws_LTW.Range("J40").Value = 2550.6667
$ws_LTW.Range("K40").Value = 1807.1666
$ws_LTW.Range("L40").Value = 2550.6667
$ws_LTW.Range("M40").Value = -1671.1666
$ws_LTW.Range("N40").Value = -2822.6667

$ws_LTW.Range("H46").Value = 2739
$ws_LTW.Range("I46").Value = 1000
$ws_LTW.Range("J46").Value = 4478
$ws_LTW.Range("K46").Value = 1000
$ws_LTW.Range("L46").Value = 4478
$ws_LTW.Range("M46").Value = -812
$ws_LTW.Range("N46").Value = -4854

$ws_LTW.Range("H122").Value = 4197.154
$ws_LTW.Range("I122").Value = 4005.7273
$ws_LTW.Range("K122").Value = 12017.1819
$ws_LTW.Range("M122").Value = -9567.1819

$ws_LTW.Range("H136").Value = 3260
$ws_LTW.Range("I136").Value = 1957.8948
$ws_LTW.Range("K136").Value = 5873.6844
$ws_LTW.Range("M136").Value = -3323.6844

$ws_LTW.Range("H140").Value = 88579.39999999999
$ws_LTW.Range("J140").Value = 88224.5
$ws_LTW.Range("L140").Value = 88224.5
$ws_LTW.Range("N140").Value = -98584.5

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H41").Value = 500001660
$ws_WVR.Range("J41").Value = 3333
$ws_WVR.Range("L41").Value = 3333
$ws_WVR.Range("N41").Value = -4113

$ws_WVR.Range("H126").Value = 55575892
$ws_WVR.Range("I126").Value = 83361336
$ws_WVR.Range("J126").Value = 5000
$ws_WVR.Range("K126").Value = 250084008
$ws_WVR.Range("L126").Value = 15000
$ws_WVR.Range("M126").Value = -250081538
$ws_WVR.Range("N126").Value = -19940

$ws_WVR.Range("H136").Value = 3383.0264
$ws_WVR.Range("I136").Value = 2206.8965
$ws_WVR.Range("J136").Value = 7172.778
$ws_WVR.Range("K136").Value = 6620.689499999999
$ws_WVR.Range("L136").Value = 21518.334
$ws_WVR.Range("M136").Value = -4070.689499999999
$ws_WVR.Range("N136").Value = -26618.334
